$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"

# Row 20: %/Uninterpretable -> b/Acknowledge (Backchannel)
$ws.Range("I20").Value = "b"
$ws.Range("J20").Value = "Acknowledge (Backchannel)"

# Row 21: sd/Statement-non-opinion -> qy/Yes-No-Question
$ws.Range("I21").Value = "qy"
$ws.Range("J21").Value = "Yes-No-Question"

# Row 24: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I24").Value = "sv"
$ws.Range("J24").Value = "Statement-opinion"

# Row 28: sd/Statement-non-opinion -> %/Uninterpretable
$ws.Range("I28").Value = "%"
$ws.Range("J28").Value = "Uninterpretable"
